# Build site at 2022-01-09 00:29:46 UTC
#
# Updates the "Ativação:" (activation) date from 01/01/2020 to 01/01/2022,
# and trims the 4th item ("Trocadores de calor tubulares" / "Tubular heat
# exchangers" and the corresponding long-form Programa/Syllabus sentences)
# out of the short/long syllabus cells, in both the Portuguese (column B)
# and English (column C) copies.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: Ativação / activation date -------------------------------
# A literal "'01/01/2022" typed straight into Excel gets auto-recognised
# as a date and reformatted, which would change the cell's number format.
# Force it in as text (leading apostrophe = quote-prefix), then copy the
# untouched format from the row below straight back on top so the cell's
# style stays exactly as it was (General, wrap/top-aligned) - only the
# text itself changes, just like the source diff.
$ws.Range("B8").Value = "'01/01/2022"
$ws.Range("B9").Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4122) | Out-Null

$ws.Range("C8").Value = "'01/01/2022"
$ws.Range("C9").Copy() | Out-Null
$ws.Range("C8").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Row 14: Programa resumido (pt) -----------------------------------
$ws.Range("B14").Value = "Perfis de temperaturas em barras de seção circular; 2) Transferência de calor por convecção; 3) Determinação do coeficiente de difusão em sistemas gás-líquido;"
$ws.Range("C14").Value = "Perfis de temperaturas em barras de seção circular; 2) Transferência de calor por convecção; 3) Determinação do coeficiente de difusão em sistemas gás-líquido;"

# --- Row 15: Short syllabus (en) ---------------------------------------
$ws.Range("B15").Value = "1) Temperature distribution in a bar with circular section; 2) Convective heat transfer; 3) Diffusion coefficient in gas-liquid systems;"
$ws.Range("C15").Value = "1) Temperature distribution in a bar with circular section; 2) Convective heat transfer; 3) Diffusion coefficient in gas-liquid systems;"

# --- Row 16: Programa (pt, long form) -----------------------------------
$ws.Range("B16").Value = "1) Perfis de temperaturas em barras de seção circular: processos envolvendo condução e convecção em barras de vários materiais e diferentes dimensões. Aplicação do princípio das aletas; 2) Transferência de calor por convecção: medidas da variação de temperatura em corpos de várias geometrias e materiais diferentes e comparação com a análise concentrada para regime transiente; 3) Determinação do coeficiente de difusão em sistemas gás-líquido: avaliação da transferência de massa entre ar e líquidos empregando tubos horizontais (célula de Stefan) em regime pseudo-estacionário;"
$ws.Range("C16").Value = "1) Perfis de temperaturas em barras de seção circular: processos envolvendo condução e convecção em barras de vários materiais e diferentes dimensões. Aplicação do princípio das aletas; 2) Transferência de calor por convecção: medidas da variação de temperatura em corpos de várias geometrias e materiais diferentes e comparação com a análise concentrada para regime transiente; 3) Determinação do coeficiente de difusão em sistemas gás-líquido: avaliação da transferência de massa entre ar e líquidos empregando tubos horizontais (célula de Stefan) em regime pseudo-estacionário;"

# --- Row 17: Syllabus (en, long form) -----------------------------------
$ws.Range("B17").Value = "1) Temperature distribution in a bar with circular section: heat transfer by conduction and convection in bars of different diameters and materials;; 2) Convective heat transfer: measures temperature variation in bodies of different geometries and materials. Comparison between the experimental data with mathematical models based on the analysis concentrated to transient parameter settings; 3) Diffusion coefficient in gas-liquid systems: analysis of mass transfer between air and liquids using horizontal pipes (Stefan cell) in pseudo-steady state;"
$ws.Range("C17").Value = "1) Temperature distribution in a bar with circular section: heat transfer by conduction and convection in bars of different diameters and materials;; 2) Convective heat transfer: measures temperature variation in bodies of different geometries and materials. Comparison between the experimental data with mathematical models based on the analysis concentrated to transient parameter settings; 3) Diffusion coefficient in gas-liquid systems: analysis of mass transfer between air and liquids using horizontal pipes (Stefan cell) in pseudo-steady state;"
